# edit.ps1
# Applies the "assets_groups_manifest" content update:
#  - renames asset groups (drop "_2" suffixes that were placeholders),
#  - rewords several descriptions,
#  - adds a new "WEB-сервера" group under "Сервисы",
#  - adds an entirely new "Сети организации" subtree (ОПС / ОСМР + 4 networks),
#  - updates sheet view (zoom 130% -> 190%, new selection) and a couple of
#    column widths that Excel recalculated as part of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so stale cells/rows from the old layout do not
# linger (the new table has more rows than the old one, and several groups
# were renamed rather than merely re-described).
$ws.Cells.Clear()

$ws.Range("A1").Value = "Название"
$ws.Range("B1").Value = "Описание"
$ws.Range("C1").Value = "Родительская группа"
$ws.Range("D1").Value = "Тип группы"
$ws.Range("E1").Value = "PDQL фильтр"
$ws.Range("F1").Value = "Доступность"
$ws.Range("G1").Value = "Косвенный ущерб"
$ws.Range("H1").Value = "Конфиденциальность"
$ws.Range("I1").Value = "Целостность"
$ws.Range("J1").Value = "Плотность целей"
$ws.Range("K1").Value = "address"
$ws.Range("L1").Value = "contactUserId"
$ws.Range("M1").Value = "internetProviders"
$ws.Range("N1").Value = "numberOfNodes"
$ws.Range("O1").Value = "registeredDomains"
$ws.Range("P1").Value = "usedNetworkApplications"
$ws.Range("Q1").Value = "usedNetworks"
$ws.Range("A2").Value = "Инфраструктура"
$ws.Range("B2").Value = "Все активы организации"
$ws.Range("C2").Value = "Root"
$ws.Range("D2").Value = "static"
$ws.Range("F2").Value = "ND"
$ws.Range("G2").Value = "ND"
$ws.Range("H2").Value = "ND"
$ws.Range("I2").Value = "ND"
$ws.Range("J2").Value = "ND"
$ws.Range("A3").Value = "Unix активы"
$ws.Range("B3").Value = "Unix подобные операционные системы"
$ws.Range("C3").Value = "Инфраструктура"
$ws.Range("D3").Value = "static"
$ws.Range("F3").Value = "ND"
$ws.Range("G3").Value = "ND"
$ws.Range("H3").Value = "ND"
$ws.Range("I3").Value = "ND"
$ws.Range("J3").Value = "ND"
$ws.Range("A4").Value = "Astra OS"
$ws.Range("B4").Value = "Хосты под управлением астра линукс"
$ws.Range("C4").Value = "Unix активы"
$ws.Range("D4").Value = "dynamic"
$ws.Range("E4").Value = "UnixHost.OsName like `"%Astra%`""
$ws.Range("F4").Value = "H"
$ws.Range("G4").Value = "ND"
$ws.Range("H4").Value = "H"
$ws.Range("I4").Value = "H"
$ws.Range("J4").Value = "ND"
$ws.Range("A5").Value = "Debian OS"
$ws.Range("B5").Value = "Хосты под управлением дебиан"
$ws.Range("C5").Value = "Unix активы"
$ws.Range("D5").Value = "dynamic"
$ws.Range("E5").Value = "UnixHost.OsName like `"%Debian%`""
$ws.Range("F5").Value = "M"
$ws.Range("G5").Value = "ND"
$ws.Range("H5").Value = "H"
$ws.Range("I5").Value = "M"
$ws.Range("J5").Value = "ND"
$ws.Range("A6").Value = "Ubuntu OS"
$ws.Range("B6").Value = "Хосты под управлением убунту"
$ws.Range("C6").Value = "Unix активы"
$ws.Range("D6").Value = "dynamic"
$ws.Range("E6").Value = "UnixHost.OsName like `"%Ubuntu%`""
$ws.Range("F6").Value = "ND"
$ws.Range("G6").Value = "ND"
$ws.Range("H6").Value = "ND"
$ws.Range("I6").Value = "ND"
$ws.Range("J6").Value = "ND"
$ws.Range("A7").Value = "Windows активы"
$ws.Range("B7").Value = "Все Windows активы организации"
$ws.Range("C7").Value = "Инфраструктура"
$ws.Range("D7").Value = "static"
$ws.Range("F7").Value = "M"
$ws.Range("G7").Value = "ND"
$ws.Range("H7").Value = "L"
$ws.Range("I7").Value = "H"
$ws.Range("J7").Value = "ND"
$ws.Range("A8").Value = "Windows АРМы"
$ws.Range("B8").Value = "Рабочие станции пользователей под управлением Windows (xp,7,8,10,11)"
$ws.Range("C8").Value = "Windows активы"
$ws.Range("D8").Value = "dynamic"
$ws.Range("E8").Value = "WindowsHost and Host.HostType = 'Desktop'"
$ws.Range("F8").Value = "L"
$ws.Range("G8").Value = "ND"
$ws.Range("H8").Value = "M"
$ws.Range("I8").Value = "L"
$ws.Range("J8").Value = "ND"
$ws.Range("A9").Value = "Windows сервера"
$ws.Range("B9").Value = "Сервера под управлением Windows (2012, 2016, 2019, 2022)"
$ws.Range("C9").Value = "Windows активы"
$ws.Range("D9").Value = "dynamic"
$ws.Range("E9").Value = "WindowsHost and Host.HostType = 'Server'"
$ws.Range("F9").Value = "H"
$ws.Range("G9").Value = "ND"
$ws.Range("H9").Value = "L"
$ws.Range("I9").Value = "H"
$ws.Range("J9").Value = "ND"
$ws.Range("A10").Value = "Сервисы"
$ws.Range("B10").Value = "Инфраструктурные сервисы организации"
$ws.Range("C10").Value = "Инфраструктура"
$ws.Range("D10").Value = "static"
$ws.Range("F10").Value = "ND"
$ws.Range("G10").Value = "ND"
$ws.Range("H10").Value = "ND"
$ws.Range("I10").Value = "ND"
$ws.Range("J10").Value = "ND"
$ws.Range("A11").Value = "DNS-сервера"
$ws.Range("B11").Value = "DNSы"
$ws.Range("C11").Value = "Сервисы"
$ws.Range("D11").Value = "dynamic"
$ws.Range("E11").Value = "Host.HostRoles.Role = 'DNS Server'"
$ws.Range("F11").Value = "M"
$ws.Range("G11").Value = "ND"
$ws.Range("H11").Value = "H"
$ws.Range("I11").Value = "H"
$ws.Range("J11").Value = "ND"
$ws.Range("A12").Value = "Контроллеры домена"
$ws.Range("B12").Value = "Контроллеры домена в организации"
$ws.Range("C12").Value = "Сервисы"
$ws.Range("D12").Value = "dynamic"
$ws.Range("E12").Value = "Host.HostRoles.Role = 'Domain Controller'"
$ws.Range("F12").Value = "H"
$ws.Range("G12").Value = "ND"
$ws.Range("H12").Value = "M"
$ws.Range("I12").Value = "H"
$ws.Range("J12").Value = "ND"
$ws.Range("A13").Value = "СУБД"
$ws.Range("B13").Value = "Хосты, на которых есть средства управления базами данных"
$ws.Range("C13").Value = "Сервисы"
$ws.Range("D13").Value = "dynamic"
$ws.Range("E13").Value = "Host.HostRoles.Role = 'Database Server'"
$ws.Range("F13").Value = "M"
$ws.Range("G13").Value = "ND"
$ws.Range("H13").Value = "L"
$ws.Range("I13").Value = "M"
$ws.Range("J13").Value = "ND"
$ws.Range("A14").Value = "WEB-сервера"
$ws.Range("B14").Value = "Хосты, на которых развернуты веб-сервисы"
$ws.Range("C14").Value = "Сервисы"
$ws.Range("D14").Value = "dynamic"
$ws.Range("E14").Value = "Host.HostRoles.Role = 'Web Server'"
$ws.Range("F14").Value = "H"
$ws.Range("G14").Value = "ND"
$ws.Range("H14").Value = "L"
$ws.Range("I14").Value = "M"
$ws.Range("J14").Value = "ND"
$ws.Range("A15").Value = "Сети организации"
$ws.Range("B15").Value = "Все сети, используемые в организации"
$ws.Range("C15").Value = "Root"
$ws.Range("D15").Value = "static"
$ws.Range("F15").Value = "ND"
$ws.Range("G15").Value = "ND"
$ws.Range("H15").Value = "ND"
$ws.Range("I15").Value = "ND"
$ws.Range("J15").Value = "ND"
$ws.Range("A16").Value = "ОПС"
$ws.Range("B16").Value = "Сети отдела прикладных систем"
$ws.Range("C16").Value = "Сети организации"
$ws.Range("D16").Value = "static"
$ws.Range("A17").Value = "ОСМР"
$ws.Range("B17").Value = "Сети отдела систем мониторинга и реагирования"
$ws.Range("C17").Value = "Сети организации"
$ws.Range("D17").Value = "static"
$ws.Range("A18").Value = "10.2.118.0/24"
$ws.Range("B18").Value = "Сеть, принадлежащая Кириллу У"
$ws.Range("C18").Value = "ОСМР"
$ws.Range("D18").Value = "dynamic"
$ws.Range("E18").Value = "Host.IpAddress in 10.2.118.0/24"
$ws.Range("F18").Value = "L"
$ws.Range("G18").Value = "ND"
$ws.Range("H18").Value = "M"
$ws.Range("I18").Value = "M"
$ws.Range("J18").Value = "ND"
$ws.Range("A19").Value = "10.2.139.0/24"
$ws.Range("B19").Value = "Сеть, принадлежащая Кириллу Б"
$ws.Range("C19").Value = "ОПС"
$ws.Range("D19").Value = "dynamic"
$ws.Range("E19").Value = "Host.IpAddress in 10.2.139.0/24"
$ws.Range("F19").Value = "M"
$ws.Range("G19").Value = "ND"
$ws.Range("H19").Value = "L"
$ws.Range("I19").Value = "H"
$ws.Range("J19").Value = "ND"
$ws.Range("A20").Value = "10.2.183.0/24"
$ws.Range("B20").Value = "Сеть, принадлежащая Кириллу К"
$ws.Range("C20").Value = "ОПС"
$ws.Range("D20").Value = "dynamic"
$ws.Range("E20").Value = "Host.IpAddress in 10.2.183.0/24"
$ws.Range("F20").Value = "H"
$ws.Range("G20").Value = "ND"
$ws.Range("H20").Value = "M"
$ws.Range("I20").Value = "H"
$ws.Range("J20").Value = "ND"
$ws.Range("A21").Value = "10.2.55.0/24"
$ws.Range("B21").Value = "Сеть, принадлежащая Виктору Т"
$ws.Range("C21").Value = "ОПС"
$ws.Range("D21").Value = "dynamic"
$ws.Range("E21").Value = "Host.IpAddress in 10.2.55.0/24"
$ws.Range("F21").Value = "H"
$ws.Range("G21").Value = "ND"
$ws.Range("H21").Value = "H"
$ws.Range("I21").Value = "H"
$ws.Range("J21").Value = "ND"

# Column formatting: B/C were widened and switched from autofit "bestFit"
# to explicit widths after the longer descriptions were entered.
$ws.Columns("B").ColumnWidth = 35.140625
$ws.Columns("C").ColumnWidth = 22.7109375

# View state captured in the commit: zoom bumped to 190% and the active
# selection left on E13.
$excel.ActiveWindow.Zoom = 190
$ws.Range("E13").Select()
